$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 785, shifting rows 785:826 down to 786:827
$ws.Rows.Item(785).Insert()

# Fill in the values for the newly inserted row 785
$ws.Cells.Item(785, 1).Value = "2026/02/07"
$ws.Cells.Item(785, 2).Value = "土"
$ws.Cells.Item(785, 3).Value = 14
$ws.Cells.Item(785, 4).Value = 201
